$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1946564885496183
$ws.Range("C2").Value = 0.5076335877862596
$ws.Range("J2").Value = 0.02290076335877863
$ws.Range("P2").Value = 0.1450381679389313
$ws.Range("S2").Value = 0.1297709923664122
$ws.Range("B3").Value = 0.01470588235294118
$ws.Range("C3").Value = 0.02205882352941177
$ws.Range("J3").Value = 0.04411764705882353
$ws.Range("P3").Value = 0.7279411764705882
$ws.Range("S3").Value = 0.1911764705882353
$ws.Range("J4").Value = 0.04444444444444445
$ws.Range("P4").Value = 0.5777777777777777
$ws.Range("S4").Value = 0.3777777777777778
$ws.Range("B6").Value = 0.07000000000000001
$ws.Range("D6").Value = 0.015
$ws.Range("F6").Value = 0.045
$ws.Range("J6").Value = 0.205
$ws.Range("O6").Value = 0.015
$ws.Range("Q6").Value = 0.195
$ws.Range("R6").Value = 0.05
$ws.Range("S6").Value = 0.405
$ws.Range("B7").Value = 0.06701030927835051
$ws.Range("D7").Value = 0.005154639175257732
$ws.Range("F7").Value = 0.04639175257731959
$ws.Range("J7").Value = 0.1237113402061856
$ws.Range("O7").Value = 0.03092783505154639
$ws.Range("Q7").Value = 0.1958762886597938
$ws.Range("R7").Value = 0.1134020618556701
$ws.Range("S7").Value = 0.4175257731958763
$ws.Range("B8").Value = 0.05597964376590331
$ws.Range("D8").Value = 0.02290076335877863
$ws.Range("F8").Value = 0.07124681933842239
$ws.Range("J8").Value = 0.09414758269720101
$ws.Range("O8").Value = 0.03307888040712468
$ws.Range("Q8").Value = 0.178117048346056
$ws.Range("R8").Value = 0.08396946564885496
$ws.Range("S8").Value = 0.460559796437659
$ws.Range("B9").Value = 0.08415841584158416
$ws.Range("D9").Value = 0.01485148514851485
$ws.Range("F9").Value = 0.07425742574257425
$ws.Range("J9").Value = 0.1138613861386139
$ws.Range("O9").Value = 0.03465346534653466
$ws.Range("Q9").Value = 0.1435643564356436
$ws.Range("R9").Value = 0.09900990099009901
$ws.Range("S9").Value = 0.4356435643564356
$ws.Range("B10").Value = 0.1104746317512275
$ws.Range("D10").Value = 0.02454991816693944
$ws.Range("F10").Value = 0.06792144026186579
$ws.Range("J10").Value = 0.1145662847790507
$ws.Range("O10").Value = 0.01227495908346972
$ws.Range("Q10").Value = 0.2176759410801964
$ws.Range("R10").Value = 0.06382978723404255
$ws.Range("S10").Value = 0.3887070376432079
$ws.Range("G11").Value = 0.1648745519713262
$ws.Range("J11").Value = 0.07526881720430108
$ws.Range("K11").Value = 0.1827956989247312
$ws.Range("L11").Value = 0.5663082437275986
$ws.Range("S11").Value = 0.01075268817204301
$ws.Range("G12").Value = 0.74375
$ws.Range("J12").Value = 0.1625
$ws.Range("K12").Value = 0.01875
$ws.Range("L12").Value = 0.04375
$ws.Range("S12").Value = 0.03125
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.2380952380952381
$ws.Range("S13").Value = 0.09523809523809523
$ws.Range("F15").Value = 0.00881057268722467
$ws.Range("H15").Value = 0.1101321585903084
$ws.Range("I15").Value = 0.09251101321585903
$ws.Range("J15").Value = 0.3612334801762114
$ws.Range("K15").Value = 0.05286343612334802
$ws.Range("M15").Value = 0.00881057268722467
$ws.Range("O15").Value = 0.1013215859030837
$ws.Range("S15").Value = 0.2643171806167401
$ws.Range("F16").Value = 0.01290322580645161
$ws.Range("H16").Value = 0.1290322580645161
$ws.Range("I16").Value = 0.07741935483870968
$ws.Range("J16").Value = 0.4967741935483871
$ws.Range("K16").Value = 0.1161290322580645
$ws.Range("M16").Value = 0.01290322580645161
$ws.Range("O16").Value = 0.05161290322580645
$ws.Range("S16").Value = 0.1032258064516129
$ws.Range("F17").Value = 0.01342281879194631
$ws.Range("H17").Value = 0.1498881431767338
$ws.Range("I17").Value = 0.06263982102908278
$ws.Range("J17").Value = 0.4250559284116331
$ws.Range("K17").Value = 0.09172259507829977
$ws.Range("M17").Value = 0.02684563758389262
$ws.Range("O17").Value = 0.07158836689038031
$ws.Range("S17").Value = 0.1588366890380313
$ws.Range("F18").Value = 0.01219512195121951
$ws.Range("H18").Value = 0.1463414634146341
$ws.Range("I18").Value = 0.1036585365853658
$ws.Range("J18").Value = 0.3780487804878049
$ws.Range("K18").Value = 0.08536585365853659
$ws.Range("M18").Value = 0.02439024390243903
$ws.Range("O18").Value = 0.08536585365853659
$ws.Range("S18").Value = 0.1646341463414634
$ws.Range("F19").Value = 0.01583710407239819
$ws.Range("H19").Value = 0.1990950226244344
$ws.Range("I19").Value = 0.09502262443438914
$ws.Range("J19").Value = 0.3687782805429864
$ws.Range("K19").Value = 0.1040723981900453
$ws.Range("M19").Value = 0.01583710407239819
$ws.Range("N19").Value = 0.003770739064856712
$ws.Range("O19").Value = 0.05882352941176471
$ws.Range("S19").Value = 0.138763197586727
